$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-insert the "被 / passive / 虛詞" dictionary row that a previous
# revert had dropped. Inserting at row 8 pushes "用掉/nv/-" and
# everything below it down by one row.
$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = "被"
$ws.Range("B8").Value = "passive"
$ws.Range("C8").Value = "虛詞"

# Call the restored row out visually: red text on a red medium box
# border drawn around A8:C8.
$rng = $ws.Range("A8:C8")
$rng.Font.Color = 255
$rng.BorderAround(1, -4138, 1, 255)

$rng.Select()
